$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1: "Correct Answer" -> "Answer"
$ws.Range("F1").Value = "Answer"

# Data validation on F1 is removed entirely (F1 drops out of the shared
# list-validation range that used to cover "F1 F3:F1048576").
$ws.Range("F1").Validation.Delete()

# The remaining list validation (now just F3:F1048576) points at the
# options row directly below it (row 3) instead of row 1.
$ws.Range("F3").Validation.Formula1 = "=B3:E3"

# Move the active selection to G12.
[void]$ws.Range("G12").Select()
